$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2 = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 0, 2.213936997104367)
    3 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987, 0, 5.582307763322248)
    4 = @(0.04172184405617529, 109.9114832445916, 18.71679738969934, 2797.565817734744, 0, 2926.235820213091)
    5 = @(0.04172184405617529, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987, 0, 1.029605918290258)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
